# "session-cookie authentication for recharge testing"
#
# The "Case2" sheet (second worksheet) lists API test cases with an
# ExpectedResponseData column (E) and an ActualResponseData header (F1)
# that was never filled in. This fills column F (rows 2-12) with the
# actual observed response for each test case - which, for this batch of
# recharge/session-cookie auth tests, matches the expected response in
# column E - completing the ActualResponseData column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case2")

$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $expected = $ws.Cells.Item($r, 5).Text
    $ws.Cells.Item($r, 6).Value = $expected
}
